# "Roboflow Annotation Report 7/28/2025" - append the new weekly progress
# row (row 67) to the bottom of the tracking table, mirroring the pattern
# of the existing rows (e.g. row 66).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write the new row's values first
$ws.Range("D67").Value = "28/7/2030"
$ws.Range("E67").Value = 380
$ws.Range("F67").Value = 950
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 1012
$ws.Range("J67").Value = "N/A"

# Copy the formatting (cell styles) from the row above so the new row
# matches the rest of the table visually
$ws.Range("D66:J66").Copy()
$ws.Range("D67:J67").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Give the new row the same height as the other data rows
$ws.Rows.Item(67).RowHeight = 15.6

# Grow the table (and its autofilter) so the new row is officially part
# of Table1
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("D4:J67"))

# Leave the selection/scroll position where the user ended up
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 57
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F74").Select() | Out-Null
